$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1636
$ws.Range("J3").Value = 1725
$ws.Range("I4").Value = 1756
$ws.Range("J4").Value = 385
$ws.Range("J6").Value = 2242
$ws.Range("I7").Value = 26202
$ws.Range("J7").Value = 6107

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 59
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 27
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 44
$ws.Range("J3").Value = 90
$ws.Range("J4").Value = 14
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 173
$ws.Range("J8").Value = 374
$ws.Range("J15").Value = 80
$ws.Range("J18").Value = 74
$ws.Range("J19").Value = 211
$ws.Range("J23").Value = 53
$ws.Range("J25").Value = 35
$ws.Range("J29").Value = 336
$ws.Range("J30").Value = 22
$ws.Range("J33").Value = 262
$ws.Range("J35").Value = 7
$ws.Range("J37").Value = 209
$ws.Range("J41").Value = 39
$ws.Range("J42").Value = 236
$ws.Range("J43").Value = 65
$ws.Range("J51").Value = 79
$ws.Range("J52").Value = 141
$ws.Range("J54").Value = 119
$ws.Range("J57").Value = 29
$ws.Range("I63").Value = 194
$ws.Range("J63").Value = 27
$ws.Range("J64").Value = 41
$ws.Range("J65").Value = 157
$ws.Range("J67").Value = 221
$ws.Range("J68").Value = 13
$ws.Range("J72").Value = 24
$ws.Range("J83").Value = 149
$ws.Range("J85").Value = 279
$ws.Range("J86").Value = 33
$ws.Range("J88").Value = 65
$ws.Range("J89").Value = 68
$ws.Range("J90").Value = 68
$ws.Range("J94").Value = 50
$ws.Range("J95").Value = 91
$ws.Range("J97").Value = 39
$ws.Range("J98").Value = 41
$ws.Range("J99").Value = 78
$ws.Range("I101").Value = 26202
$ws.Range("J101").Value = 6107

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 33
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 262

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 100
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 336

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 52
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 211

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 70
$ws.Range("J3").Value = 110
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 279

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 49
$ws.Range("J6").Value = 126
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 46
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 141

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J2").Value = 11
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 7

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 122
$ws.Range("J3").Value = 128
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 374

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 56
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 173
